$wb = $excel.ActiveWorkbook

$er = $wb.Worksheets.Item("ER")
$erList = $wb.Worksheets.Item("ERList")

# --- Add two new trailing sheets: "Sheet1" and "Sheet2" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet1 = $wb.Worksheets.Add($null, $lastSheet)
$sheet1.Name = "Sheet1"
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet2.Name = "Sheet2"

# --- "Sheet1" gets a snapshot of ERList's two data rows (as they were) ---
$null = $erList.Rows("2:3").Copy()
$null = $sheet1.Rows("2:3").Select()
$null = $sheet1.Paste()

$null = $sheet1.Hyperlinks.Add($sheet1.Range("B2"), "mailto:externalreviewer5@intrees.org")
$null = $sheet1.Hyperlinks.Add($sheet1.Range("B3"), "mailto:externalreviewer4@intrees.org")
$null = $sheet1.Hyperlinks.Add($sheet1.Range("C2"), "mailto:Test@123")
$null = $sheet1.Hyperlinks.Add($sheet1.Range("C3"), "mailto:Test@123")
$sheet1.Range("B2:C3").Style = "Hyperlink"

# --- Update ERList with the new reviewer names / emails ---
$erList.Range("A2").Value = "Alyazia Khamis"
$erList.Range("B2").Value = "ertesting2he@gmail.com"
$erList.Range("A3").Value = "Owaisha Aamer"
$erList.Range("B3").Value = "er3hadeel@gmail.com"

# recreate the hyperlinks on ERList against the new data
$null = $erList.Range("A1").Hyperlinks.Delete()
$null = $erList.Hyperlinks.Add($erList.Range("B3"), "mailto:er3hadeel@gmail.com")
$null = $erList.Hyperlinks.Add($erList.Range("B2"), "mailto:ertesting2he@gmail.com")
$null = $erList.Hyperlinks.Add($erList.Range("C2"), "mailto:Test@123")
$null = $erList.Hyperlinks.Add($erList.Range("C3"), "mailto:Test@123")
$erList.Range("B2:C3").Style = "Hyperlink"
$null = $erList.Range("C3").Select()

# --- Make ER the active/selected tab (was RequestToChangeER before) ---
$null = $er.Select()
$null = $er.Range("A2").Select()
